# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.805.59"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "3.143.63"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'592.28"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'145.45"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.135.79"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").Value = "'5.88"
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("E12").Value = "  -2.18%  "
$ws.Range("E13").Value = "  -3.10%  "
$ws.Range("D14").Value = "'37.19"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").Value = "3.666.66"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").Value = "3.141.19"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "63.699.87"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").Value = "'468.05"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "'14.35"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("D24").Value = "'12.98"
$ws.Range("E24").Value = "  -2.46%  "
$ws.Range("D25").Value = "'81.34"
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("E26").Value = "  +5.62%  "
$ws.Range("D28").Value = "'9.77"
$ws.Range("E28").Value = "  +8.55%  "
$ws.Range("D29").Value = "'7.38"
$ws.Range("E29").Value = "  +7.35%  "
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").Value = "'27.74"
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("D35").Value = "0.0₃0843"
$ws.Range("E35").Value = "  -4.12%  "
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "'2.31"
$ws.Range("E37").Value = "  -3.51%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "'6.15"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("E39").Value = "  -5.44%  "
$ws.Range("D40").Value = "'51.40"
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "'9.32"
$ws.Range("E41").Value = "  +6.63%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'455.76"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("E43").Value = "  +5.03%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0371"
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.915.28"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").Value = "'39.70"
$ws.Range("E46").Value = "  +11.20%  "
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("D48").Value = "'130.18"
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("D50").Value = "'2.24"
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("E51").Value = "  -1.28%  "
